$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.590.62'
$ws.Range('E2').Value = '  +1.02%  '
$ws.Range('D3').Value = '3.318.38'
$ws.Range('E3').Value = '  +5.56%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.79'
$ws.Range('E5').Value = '  +3.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.97'
$ws.Range('E6').Value = '  +3.35%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.316.43'
$ws.Range('E8').Value = '  +5.57%  '
$ws.Range('E9').Value = '  +0.89%  '
$ws.Range('E10').Value = '  +3.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.50'
$ws.Range('E11').Value = '  +4.15%  '
$ws.Range('E12').Value = '  +2.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000247'
$ws.Range('E13').Value = '  +1.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.86'
$ws.Range('E14').Value = '  +2.45%  '
$ws.Range('D15').Value = '3.864.95'
$ws.Range('E15').Value = '  +5.56%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.121'
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('D17').Value = '3.320.44'
$ws.Range('E17').Value = '  +5.52%  '
$ws.Range('D18').Value = '63.663.63'
$ws.Range('E18').Value = '  +1.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.85'
$ws.Range('E19').Value = '  +3.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '479.40'
$ws.Range('E20').Value = '  +1.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.11'
$ws.Range('E21').Value = '  +1.19%  '
$ws.Range('E22').Value = '  +4.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.17'
$ws.Range('E23').Value = '  +6.12%  '
$ws.Range('E24').Value = '  +6.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.07'
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  +2.73%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.22'
$ws.Range('E29').Value = '  +2.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.18'
$ws.Range('E30').Value = '  +3.53%  '
$ws.Range('E31').Value = '  +2.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '29.01'
$ws.Range('E32').Value = '  +8.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.106'
$ws.Range('E33').Value = '  +1.61%  '
$ws.Range('E34').Value = '  +0.41%  '
$ws.Range('E35').Value = '  +4.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.04'
$ws.Range('E36').Value = '  +4.84%  '
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = '0.0₃0744'
$ws.Range('E37').Value = '  +7.15%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '51.99'
$ws.Range('E38').Value = '  -0.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0399'
$ws.Range('E39').Value = '  +3.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '433.72'
$ws.Range('E40').Value = '  +4.47%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.122'
$ws.Range('E41').Value = '  +10.47%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '3.084.56'
$ws.Range('E42').Value = '  +5.73%  '
$ws.Range('E43').Value = '  +0.99%  '
$ws.Range('E44').Value = '  +0.98%  '
$ws.Range('E45').Value = '  +2.48%  '
$ws.Range('E46').Value = '  +4.73%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '36.97'
$ws.Range('E47').Value = '  +14.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '26.26'
$ws.Range('E48').Value = '  +3.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.999'
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('E50').Value = '  +1.09%  '
$ws.Range('E51').Value = '  +2.89%  '
